# Update UI class diagram in Developer Guide
#
# 1) Rename three shapes on the (single) slide to reflect the renamed
#    UI components (BrowserPanel -> TaskCardHeader, PersonListPanel ->
#    TaskListPanel, PersonCard -> TaskCard).
# 2) Refresh the cached "datetimeFigureOut" field text (1/7/2017 ->
#    3/27/2017) wherever it appears in the slide master and the slide
#    layouts.

$p = $ppt.ActivePresentation

# --- 1) Rename the renamed UI-component shapes on slide 1 -----------------
$s = $p.Slides.Item(1)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $t = $shp.TextFrame.TextRange.Text
        if ($t -eq "BrowserPanel") {
            $shp.TextFrame.TextRange.Text = "TaskCardHeader"
        } elseif ($t -eq "PersonListPanel") {
            $shp.TextFrame.TextRange.Text = "TaskListPanel"
        } elseif ($t -eq "PersonCard") {
            $shp.TextFrame.TextRange.Text = "TaskCard"
        }
    }
}

# --- 2) Re-stamp the cached date field text --------------------------------
$oldDate = "1/7/2017"
$newDate = "3/27/2017"

$master = $p.SlideMaster

# Slide master's own "Date Placeholder" shape.
for ($j = 1; $j -le $master.Shapes.Count; $j++) {
    $shp = $master.Shapes.Item($j)
    if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -eq $oldDate) {
        $shp.TextFrame.TextRange.Text = $newDate
    }
}

# Every slide layout's "Date Placeholder" shape.
$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $shp = $layout.Shapes.Item($j)
        if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -eq $oldDate) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}
